$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.800.67"
$ws.Range("E2").Value = "  +4.25%  "
$ws.Range("D3").Value = "3.354.20"
$ws.Range("E3").Value = "  +4.36%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +3.63%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "3.931.83"
$ws.Range("E12").Value = "  +4.29%  "
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "62.823.11"
$ws.Range("E16").Value = "  +4.20%  "
$ws.Range("D17").Value = "3.355.80"
$ws.Range("E17").Value = "  +4.10%  "
$ws.Range("E18").Value = "  +0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.60%  "
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  +5.28%  "
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "0.0₃0952"
$ws.Range("E27").Value = "  +4.68%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.31%  "
$ws.Range("E30").Value = "  +3.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("E36").Value = "  +8.03%  "
$ws.Range("E37").Value = "  +12.02%  "
$ws.Range("E38").Value = "  +3.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0743"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.90%  "
$ws.Range("D40").Value = "2.811.62"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  +6.75%  "
$ws.Range("E42").Value = "  +3.56%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("D46").Value = "3.398.41"
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("E47").Value = "  +5.26%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "288.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.60%  "
$ws.Range("E51").Value = "  -1.43%  "
